$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '22.027.72'
$c.NumberFormat = "General"
$c.ClearFormats()
$ws.Range("E2").Value = '  -1.94%  '

# Row 3
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '1.555.08'
$c.NumberFormat = "General"
$c.ClearFormats()
$ws.Range("E3").Value = '  -1.21%  '

# Row 4
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.002'
$c.NumberFormat = "General"
$c.ClearFormats()
$ws.Range("E4").Value = '  +0.09%  '

# Row 5
$ws.Range("E5").Value = '  +0.12%  '

# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '286.48'
$c.NumberFormat = "General"
$c.ClearFormats()
$ws.Range("E6").Value = '  -0.49%  '

# Row 7
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.3745'
$c.NumberFormat = "General"
$c.ClearFormats()
$ws.Range("E7").Value = '  +0.89%  '

# Row 8
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.3230'
$c.NumberFormat = "General"
$c.ClearFormats()
$ws.Range("E8").Value = '  -2.97%  '

# Row 9
$ws.Range("B9").Value = 'OKB'
$ws.Range("C9").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '41.51'
$c.NumberFormat = "General"
$c.ClearFormats()
$ws.Range("E9").Value = '  -12.70%  '

# Row 10
$ws.Range("B10").Value = 'Polygon'
$ws.Range("C10").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '1.128'
$c.NumberFormat = "General"
$c.ClearFormats()
$ws.Range("E10").Value = '  -1.90%  '

# Row 11
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.07300'
$c.NumberFormat = "General"
$c.ClearFormats()
$ws.Range("E11").Value = '  -3.42%  '

# Row 12
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '1.002'
$c.NumberFormat = "General"
$c.ClearFormats()
$ws.Range("E12").Value = '  +0.10%  '

# Row 13
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '19.71'
$c.NumberFormat = "General"
$c.ClearFormats()

# Row 14
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '6.852'
$c.NumberFormat = "General"
$c.ClearFormats()
$ws.Range("E14").Value = '  -1.45%  '

# Row 15
$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '5.687'
$c.NumberFormat = "General"
$c.ClearFormats()
$ws.Range("E15").Value = '  -4.49%  '

# Row 16
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '1.562.64'
$c.NumberFormat = "General"
$c.ClearFormats()
$ws.Range("E16").Value = '  -0.29%  '

# Row 17
$ws.Range("E17").Value = '  -3.57%  '

# Row 18
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '0.06624'
$c.NumberFormat = "General"
$c.ClearFormats()
$ws.Range("E18").Value = '  -1.62%  '

# Row 19
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '85.02'
$c.NumberFormat = "General"
$c.ClearFormats()
$ws.Range("E19").Value = '  -3.66%  '

# Row 20
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '6.449'
$c.NumberFormat = "General"
$c.ClearFormats()
$ws.Range("E20").Value = '  +0.65%  '

# Row 21
$ws.Range("E21").Value = '  +0.12%  '

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '15.97'
$c.NumberFormat = "General"
$c.ClearFormats()
$ws.Range("E22").Value = '  -3.48%  '

# Row 23
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '11.56'
$c.NumberFormat = "General"
$c.ClearFormats()
$ws.Range("E23").Value = '  -3.89%  '

# Row 24
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '22.039.54'
$c.NumberFormat = "General"
$c.ClearFormats()
$ws.Range("E24").Value = '  -1.84%  '

# Row 25
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '2.261'
$c.NumberFormat = "General"
$c.ClearFormats()
$ws.Range("E25").Value = '  -5.50%  '

# Row 26
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '2.518'
$c.NumberFormat = "General"
$c.ClearFormats()
$ws.Range("E26").Value = '  -4.21%  '

# Row 27
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '149.51'
$c.NumberFormat = "General"
$c.ClearFormats()
$ws.Range("E27").Value = '  -1.26%  '

# Row 28
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '18.84'
$c.NumberFormat = "General"
$c.ClearFormats()
$ws.Range("E28").Value = '  -4.24%  '

# Row 29
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '4.845'
$c.NumberFormat = "General"
$c.ClearFormats()
$ws.Range("E29").Value = '  -2.59%  '

# Row 30
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '1.740.33'
$c.NumberFormat = "General"
$c.ClearFormats()
$ws.Range("E30").Value = '  -1.28%  '

# Row 31
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '120.10'
$c.NumberFormat = "General"
$c.ClearFormats()
$ws.Range("E31").Value = '  -4.24%  '

# Row 32
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '1.110'
$c.NumberFormat = "General"
$c.ClearFormats()
$ws.Range("E32").Value = '  +1.52%  '

# Row 33
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '5.956'
$c.NumberFormat = "General"
$c.ClearFormats()
$ws.Range("E33").Value = '  -2.47%  '

# Row 34
$ws.Range("B34").Value = 'FraxShare'
$ws.Range("C34").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '9.269'
$c.NumberFormat = "General"
$c.ClearFormats()
$ws.Range("E34").Value = '  -6.25%  '

# Row 35
$ws.Range("B35").Value = 'Stellar'
$ws.Range("C35").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.08114'
$c.NumberFormat = "General"
$c.ClearFormats()
$ws.Range("E35").Value = '  -2.88%  '

# Row 36
$ws.Range("B36").Value = 'WEMIXTOKEN'
$ws.Range("C36").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '1.583'
$c.NumberFormat = "General"
$c.ClearFormats()
$ws.Range("E36").Value = '  -20.48%  '

# Row 37
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '5.226'
$c.NumberFormat = "General"
$c.ClearFormats()
$ws.Range("E37").Value = '  -2.47%  '

# Row 38
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.02272'
$c.NumberFormat = "General"
$c.ClearFormats()
$ws.Range("E38").Value = '  -7.87%  '

# Row 39
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.06126'
$c.NumberFormat = "General"
$c.ClearFormats()
$ws.Range("E39").Value = '  -4.06%  '

# Row 40
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.2118'
$c.NumberFormat = "General"
$c.ClearFormats()
$ws.Range("E40").Value = '  -5.38%  '

# Row 41
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '1.210'
$c.NumberFormat = "General"
$c.ClearFormats()
$ws.Range("E41").Value = '  -7.37%  '

# Row 42
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '10.90'
$c.NumberFormat = "General"
$c.ClearFormats()
$ws.Range("E42").Value = '  -4.96%  '

# Row 43
$ws.Range("E43").Value = '  +0.18%  '

# Row 44
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.5933'
$c.NumberFormat = "General"
$c.ClearFormats()
$ws.Range("E44").Value = '  -5.46%  '

# Row 45
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '13.60'
$c.NumberFormat = "General"
$c.ClearFormats()
$ws.Range("E45").Value = '  -3.35%  '

# Row 46
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '3.719'
$c.NumberFormat = "General"
$c.ClearFormats()
$ws.Range("E46").Value = '  -1.58%  '

# Row 47
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.5736'
$c.NumberFormat = "General"
$c.ClearFormats()
$ws.Range("E47").Value = '  -5.98%  '

# Row 48
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '1.942'
$c.NumberFormat = "General"
$c.ClearFormats()
$ws.Range("E48").Value = '  -5.44%  '

# Row 49
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '119.64'
$c.NumberFormat = "General"
$c.ClearFormats()
$ws.Range("E49").Value = '  -4.54%  '

# Row 50
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '1.155'
$c.NumberFormat = "General"
$c.ClearFormats()

# Row 51
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.06932'
$c.NumberFormat = "General"
$c.ClearFormats()
$ws.Range("E51").Value = '  -3.96%  '
